$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 45 (pushes the existing data for rows 45-172 down
# to rows 47-174, matching the growth of the sheet's used range from R172 to R174).
$ws.Rows.Item(45).Resize(2).Insert()

# New row 45: Feria Lagunitas de Puerto Montt, Apio, Americana (o), Primera
$ws.Cells.Item(45, 1).Value = 4
$ws.Cells.Item(45, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(45, 3).Value = "Los Lagos"
$ws.Cells.Item(45, 4).Value = 44544
$ws.Cells.Item(45, 5).Value = 10
$ws.Cells.Item(45, 6).Value = 100112017
$ws.Cells.Item(45, 7).Value = "Apio"
$ws.Cells.Item(45, 8).Value = "Americana (o)"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 20
$ws.Cells.Item(45, 11).Value = 12000
$ws.Cells.Item(45, 12).Value = 12000
$ws.Cells.Item(45, 13).Value = 12000
$ws.Cells.Item(45, 14).Value = '$/docena de matas'
$ws.Cells.Item(45, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(45, 16).Value = 2000
$ws.Cells.Item(45, 17).Value = 6
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# New row 46: Feria Lagunitas de Puerto Montt, Apio, Americana (o), Segunda
$ws.Cells.Item(46, 1).Value = 4
$ws.Cells.Item(46, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(46, 3).Value = "Los Lagos"
$ws.Cells.Item(46, 4).Value = 44544
$ws.Cells.Item(46, 5).Value = 10
$ws.Cells.Item(46, 6).Value = 100112017
$ws.Cells.Item(46, 7).Value = "Apio"
$ws.Cells.Item(46, 8).Value = "Americana (o)"
$ws.Cells.Item(46, 9).Value = "Segunda"
$ws.Cells.Item(46, 10).Value = 20
$ws.Cells.Item(46, 11).Value = 10000
$ws.Cells.Item(46, 12).Value = 10000
$ws.Cells.Item(46, 13).Value = 10000
$ws.Cells.Item(46, 14).Value = '$/docena de matas'
$ws.Cells.Item(46, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(46, 16).Value = 1667
$ws.Cells.Item(46, 17).Value = 6
$ws.Cells.Item(46, 18).Value = "Hortaliza"
